$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data to append: row, D(survey_round), E(panel letter), F(wave, may be text for split waves), G(date serial), H(spss_name)
$rows = @(
    @{R=71; D=59; E="F"; F=20;    G=44329; H="21-037554_PFW20_Final_ICUO"},
    @{R=72; D=60; E="E"; F=21;    G=44337; H="21-037558_PEW21_Final_ICUO"},
    @{R=73; D=61; E="F"; F=21;    G=44344; H="21-037554_PFW21_Final_icuo"},
    @{R=74; D=62; E="E"; F=22;    G=44351; H="21-037558_PEW22_Final_ICUO"},
    @{R=75; D=63; E="F"; F=22;    G=44358; H="21-037554_PFW22_Final_icuo"},
    @{R=76; D=64; E="E"; F=23;    G=44365; H="21-037558_PEW23_Final_ICUO"},
    @{R=77; D=65; E="F"; F="23a"; G=44372; H="21-037554_PFW23a_Final_ICUO"},
    @{R=78; D=65; E="F"; F="23b"; G=44377; H="21-037554_PFW23b_Final_ICUO"},
    @{R=79; D=66; E="E"; F=24;    G=44379; H="21-037558_PEW24_Final_ICUO"},
    @{R=80; D=67; E="F"; F=24;    G=44389; H="21-037554_PFW24_Final_ICUO"},
    @{R=81; D=68; E="E"; F=25;    G=44392; H="21-037558_PEW25_Final_ICUO"},
    @{R=82; D=69; E="F"; F=25;    G=44399; H="21-037554_PFW25_Final_ICUO"}
)

foreach ($row in $rows) {
    $r = $row.R
    $prev = $r - 1

    # Copy formatting (number format etc.) from the row above so new cells match existing style (e.g. date format on column G)
    $ws.Range("A" + $prev + ":J" + $prev).Copy()
    $ws.Range("A" + $r + ":J" + $r).PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = 3
    $ws.Cells.Item($r, 2).Value = 0
    $ws.Cells.Item($r, 3).Value = "uk"
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 7).Value = $row.G
    # H (spss_name) then F (wave) - matches original authoring/entry order
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 9).Formula = "=C$r&""_""&""sr""&TEXT(D$r,""00"")&""_""&YEAR(G$r)&TEXT(G$r,""MM"")&TEXT(G$r,""DD"")&""_p""&E$r&""_wv""&TEXT(F$r,""00"")&"""""
    $ws.Cells.Item($r, 10).Value = 1
}

# Reflect the new scroll position / active selection used when the data was last reviewed
$win = $excel.ActiveWindow
$win.ScrollRow = 63
$win.ScrollColumn = 1
$ws.Range("J82").Select() | Out-Null
